$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first (leftmost) column of the sheet, containing numeric internal
# IDs, is being removed entirely. This shifts the former "product name"
# column (B) to A, "gtin" (C) to B, and "ntin" (D) to C, along with all
# of their data/formulas and column widths.
$ws.Columns("A:A").Select()
$ws.Columns("A:A").Delete()
